$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I; this shifts the existing
# "Obsolescence percentage" column (and its data/formatting) from I to J.
$ws.Columns.Item(9).Insert()

# New column header: "Distribution channel code" (bold, matching the
# other header cells in row 1).
$ws.Cells.Item(1, 9).Value = "Distribution channel code"
$ws.Cells.Item(1, 9).Font.Bold = $true

# New column data.
$ws.Cells.Item(2, 9).Value = "TR"
$ws.Cells.Item(3, 9).Value = "GO"

# New column's width (close to the author's slightly-adjusted custom width).
$ws.Columns.Item(9).ColumnWidth = 21.73
